$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Utah Jazz" -------------------------------------------------
# Consolidate the five separate RK/PLAYER/STAT mini-tables (cols A-V, rows 1-11)
# plus the second block of four mini-tables (cols A-Q, rows 13-23) into one
# single wide table (cols A-R, rows 1-11).
$ws1 = $wb.Worksheets.Item("Utah Jazz")

# Wipe the whole old layout (content + formatting) before rebuilding it.
$ws1.Range("A1:V23").Clear()

# Header row
$ws1.Cells.Item(1, 1).Value2 = "PLAYER"
$ws1.Cells.Item(1, 2).Value2 = "Games"
$ws1.Cells.Item(1, 3).Value2 = "PLAYER1"
$ws1.Cells.Item(1, 4).Value2 = "MinutesPlayed"
$ws1.Cells.Item(1, 5).Value2 = "PLAYER2"
$ws1.Cells.Item(1, 6).Value2 = "FieldGoals"
$ws1.Cells.Item(1, 7).Value2 = "PLAYER3"
$ws1.Cells.Item(1, 8).Value2 = "PtFieldGoals"
$ws1.Cells.Item(1, 9).Value2 = "PLAYER4"
$ws1.Cells.Item(1, 10).Value2 = "FreeThrows"
$ws1.Cells.Item(1, 11).Value2 = "PLAYER5"
$ws1.Cells.Item(1, 12).Value2 = "TotalRebounds"
$ws1.Cells.Item(1, 13).Value2 = "PLAYER6"
$ws1.Cells.Item(1, 14).Value2 = "Assists"
$ws1.Cells.Item(1, 15).Value2 = "PLAYER7"
$ws1.Cells.Item(1, 16).Value2 = "Steals"
$ws1.Cells.Item(1, 17).Value2 = "PLAYER8"
$ws1.Cells.Item(1, 18).Value2 = "Blocks"

# The three "pair-total" headers (MinutesPlayed, FieldGoals, PtFieldGoals) keep
# the centered style the stat column already carried in the source mini-tables.
$ws1.Cells.Item(1, 4).HorizontalAlignment = -4108
$ws1.Cells.Item(1, 6).HorizontalAlignment = -4108
$ws1.Cells.Item(1, 8).HorizontalAlignment = -4108

# Data rows 2-11
# Row 2
$ws1.Cells.Item(2, 1).Value2 = "John Stockton*"
$ws1.Cells.Item(2, 2).Value2 = 1504
$ws1.Cells.Item(2, 3).Value2 = "Karl Malone*"
$ws1.Cells.Item(2, 4).Value2 = 53479
$ws1.Cells.Item(2, 5).Value2 = "Karl Malone*"
$ws1.Cells.Item(2, 6).Value2 = 13335
$ws1.Cells.Item(2, 7).Value2 = "Joe Ingles"
$ws1.Cells.Item(2, 8).Value2 = 1071
$ws1.Cells.Item(2, 9).Value2 = "Karl Malone*"
$ws1.Cells.Item(2, 10).Value2 = 9619
$ws1.Cells.Item(2, 11).Value2 = "Karl Malone*"
$ws1.Cells.Item(2, 12).Value2 = 14601
$ws1.Cells.Item(2, 13).Value2 = "John Stockton*"
$ws1.Cells.Item(2, 14).Value2 = 15806
$ws1.Cells.Item(2, 15).Value2 = "John Stockton*"
$ws1.Cells.Item(2, 16).Value2 = 3265
$ws1.Cells.Item(2, 17).Value2 = "Mark Eaton"
$ws1.Cells.Item(2, 18).Value2 = 3064

# Row 3
$ws1.Cells.Item(3, 1).Value2 = "Karl Malone*"
$ws1.Cells.Item(3, 2).Value2 = 1434
$ws1.Cells.Item(3, 3).Value2 = "John Stockton*"
$ws1.Cells.Item(3, 4).Value2 = 47764
$ws1.Cells.Item(3, 5).Value2 = "John Stockton*"
$ws1.Cells.Item(3, 6).Value2 = 7039
$ws1.Cells.Item(3, 7).Value2 = "Donovan Mitchell"
$ws1.Cells.Item(3, 8).Value2 = 958
$ws1.Cells.Item(3, 9).Value2 = "John Stockton*"
$ws1.Cells.Item(3, 10).Value2 = 4788
$ws1.Cells.Item(3, 11).Value2 = "Rudy Gobert"
$ws1.Cells.Item(3, 12).Value2 = 7119
$ws1.Cells.Item(3, 13).Value2 = "Karl Malone*"
$ws1.Cells.Item(3, 14).Value2 = 5085
$ws1.Cells.Item(3, 15).Value2 = "Karl Malone*"
$ws1.Cells.Item(3, 16).Value2 = 2035
$ws1.Cells.Item(3, 17).Value2 = "Andrei Kirilenko"
$ws1.Cells.Item(3, 18).Value2 = 1380

# Row 4
$ws1.Cells.Item(4, 1).Value2 = "Mark Eaton"
$ws1.Cells.Item(4, 2).Value2 = 875
$ws1.Cells.Item(4, 3).Value2 = "Mark Eaton"
$ws1.Cells.Item(4, 4).Value2 = 25169
$ws1.Cells.Item(4, 5).Value2 = "Darrell Griffith"
$ws1.Cells.Item(4, 6).Value2 = 5237
$ws1.Cells.Item(4, 7).Value2 = "John Stockton*"
$ws1.Cells.Item(4, 8).Value2 = 845
$ws1.Cells.Item(4, 9).Value2 = "Adrian Dantley*"
$ws1.Cells.Item(4, 10).Value2 = 3814
$ws1.Cells.Item(4, 11).Value2 = "Mark Eaton"
$ws1.Cells.Item(4, 12).Value2 = 6939
$ws1.Cells.Item(4, 13).Value2 = "Rickey Green"
$ws1.Cells.Item(4, 14).Value2 = 4159
$ws1.Cells.Item(4, 15).Value2 = "Rickey Green"
$ws1.Cells.Item(4, 16).Value2 = 1100
$ws1.Cells.Item(4, 17).Value2 = "Rudy Gobert"
$ws1.Cells.Item(4, 18).Value2 = 1357

# Row 5
$ws1.Cells.Item(5, 1).Value2 = "Darrell Griffith"
$ws1.Cells.Item(5, 2).Value2 = 765
$ws1.Cells.Item(5, 3).Value2 = "Darrell Griffith"
$ws1.Cells.Item(5, 4).Value2 = 21403
$ws1.Cells.Item(5, 5).Value2 = "Adrian Dantley*"
$ws1.Cells.Item(5, 6).Value2 = 4908
$ws1.Cells.Item(5, 7).Value2 = "Gordon Hayward"
$ws1.Cells.Item(5, 8).Value2 = 689
$ws1.Cells.Item(5, 9).Value2 = "Andrei Kirilenko"
$ws1.Cells.Item(5, 10).Value2 = 2520
$ws1.Cells.Item(5, 11).Value2 = "Derrick Favors"
$ws1.Cells.Item(5, 12).Value2 = 4626
$ws1.Cells.Item(5, 13).Value2 = "Deron Williams"
$ws1.Cells.Item(5, 14).Value2 = 4003
$ws1.Cells.Item(5, 15).Value2 = "Andrei Kirilenko"
$ws1.Cells.Item(5, 16).Value2 = 960
$ws1.Cells.Item(5, 17).Value2 = "Greg Ostertag"
$ws1.Cells.Item(5, 18).Value2 = 1253

# Row 6
$ws1.Cells.Item(6, 1).Value2 = "Thurl Bailey"
$ws1.Cells.Item(6, 2).Value2 = 708
$ws1.Cells.Item(6, 3).Value2 = "Andrei Kirilenko"
$ws1.Cells.Item(6, 4).Value2 = 20989
$ws1.Cells.Item(6, 5).Value2 = "Thurl Bailey"
$ws1.Cells.Item(6, 6).Value2 = 3989
$ws1.Cells.Item(6, 7).Value2 = "Jordan Clarkson"
$ws1.Cells.Item(6, 8).Value2 = 574
$ws1.Cells.Item(6, 9).Value2 = "Rudy Gobert"
$ws1.Cells.Item(6, 10).Value2 = 2000
$ws1.Cells.Item(6, 11).Value2 = "John Stockton*"
$ws1.Cells.Item(6, 12).Value2 = 4051
$ws1.Cells.Item(6, 13).Value2 = "Joe Ingles"
$ws1.Cells.Item(6, 14).Value2 = 2213
$ws1.Cells.Item(6, 15).Value2 = "Darrell Griffith"
$ws1.Cells.Item(6, 16).Value2 = 931
$ws1.Cells.Item(6, 17).Value2 = "Karl Malone*"
$ws1.Cells.Item(6, 18).Value2 = 1125

# Row 7
$ws1.Cells.Item(7, 1).Value2 = "Greg Ostertag"
$ws1.Cells.Item(7, 2).Value2 = 700
$ws1.Cells.Item(7, 3).Value2 = "Thurl Bailey"
$ws1.Cells.Item(7, 4).Value2 = 20523
$ws1.Cells.Item(7, 5).Value2 = "Pete Maravich*"
$ws1.Cells.Item(7, 6).Value2 = 3258
$ws1.Cells.Item(7, 7).Value2 = "Bojan Bogdanović"
$ws1.Cells.Item(7, 8).Value2 = 550
$ws1.Cells.Item(7, 9).Value2 = "Gordon Hayward"
$ws1.Cells.Item(7, 10).Value2 = 1946
$ws1.Cells.Item(7, 11).Value2 = "Greg Ostertag"
$ws1.Cells.Item(7, 12).Value2 = 3978
$ws1.Cells.Item(7, 13).Value2 = "Andrei Kirilenko"
$ws1.Cells.Item(7, 14).Value2 = 1919
$ws1.Cells.Item(7, 15).Value2 = "Bryon Russell"
$ws1.Cells.Item(7, 16).Value2 = 728
$ws1.Cells.Item(7, 17).Value2 = "Thurl Bailey"
$ws1.Cells.Item(7, 18).Value2 = 879

# Row 8
$ws1.Cells.Item(8, 1).Value2 = "Andrei Kirilenko"
$ws1.Cells.Item(8, 2).Value2 = 681
$ws1.Cells.Item(8, 3).Value2 = "Rudy Gobert"
$ws1.Cells.Item(8, 4).Value2 = 18301
$ws1.Cells.Item(8, 5).Value2 = "Derrick Favors"
$ws1.Cells.Item(8, 6).Value2 = 2976
$ws1.Cells.Item(8, 7).Value2 = "Bryon Russell"
$ws1.Cells.Item(8, 8).Value2 = 540
$ws1.Cells.Item(8, 9).Value2 = "Thurl Bailey"
$ws1.Cells.Item(8, 10).Value2 = 1915
$ws1.Cells.Item(8, 11).Value2 = "Rich Kelley"
$ws1.Cells.Item(8, 12).Value2 = 3972
$ws1.Cells.Item(8, 13).Value2 = "Jeff Hornacek"
$ws1.Cells.Item(8, 14).Value2 = 1895
$ws1.Cells.Item(8, 15).Value2 = "Jeff Hornacek"
$ws1.Cells.Item(8, 16).Value2 = 618
$ws1.Cells.Item(8, 17).Value2 = "Derrick Favors"
$ws1.Cells.Item(8, 18).Value2 = 840

# Row 9
$ws1.Cells.Item(9, 1).Value2 = "Derrick Favors"
$ws1.Cells.Item(9, 2).Value2 = 644
$ws1.Cells.Item(9, 3).Value2 = "Adrian Dantley*"
$ws1.Cells.Item(9, 4).Value2 = 17899
$ws1.Cells.Item(9, 5).Value2 = "Donovan Mitchell"
$ws1.Cells.Item(9, 6).Value2 = 2953
$ws1.Cells.Item(9, 7).Value2 = "Darrell Griffith"
$ws1.Cells.Item(9, 8).Value2 = 530
$ws1.Cells.Item(9, 9).Value2 = "Pete Maravich*"
$ws1.Cells.Item(9, 10).Value2 = 1801
$ws1.Cells.Item(9, 11).Value2 = "Thurl Bailey"
$ws1.Cells.Item(9, 12).Value2 = 3881
$ws1.Cells.Item(9, 13).Value2 = "Pete Maravich*"
$ws1.Cells.Item(9, 14).Value2 = 1844
$ws1.Cells.Item(9, 15).Value2 = "Paul Millsap"
$ws1.Cells.Item(9, 16).Value2 = 604
$ws1.Cells.Item(9, 17).Value2 = "Paul Millsap"
$ws1.Cells.Item(9, 18).Value2 = 520

# Row 10
$ws1.Cells.Item(10, 1).Value2 = "Bryon Russell"
$ws1.Cells.Item(10, 2).Value2 = 628
$ws1.Cells.Item(10, 3).Value2 = "Rickey Green"
$ws1.Cells.Item(10, 4).Value2 = 17329
$ws1.Cells.Item(10, 5).Value2 = "Carlos Boozer"
$ws1.Cells.Item(10, 6).Value2 = 2804
$ws1.Cells.Item(10, 7).Value2 = "Mehmet Okur"
$ws1.Cells.Item(10, 8).Value2 = 517
$ws1.Cells.Item(10, 9).Value2 = "Mehmet Okur"
$ws1.Cells.Item(10, 10).Value2 = 1648
$ws1.Cells.Item(10, 11).Value2 = "Andrei Kirilenko"
$ws1.Cells.Item(10, 12).Value2 = 3836
$ws1.Cells.Item(10, 13).Value2 = "Gordon Hayward"
$ws1.Cells.Item(10, 14).Value2 = 1762
$ws1.Cells.Item(10, 15).Value2 = "Joe Ingles"
$ws1.Cells.Item(10, 16).Value2 = 544
$ws1.Cells.Item(10, 17).Value2 = "Ben Poquette"
$ws1.Cells.Item(10, 18).Value2 = 517

# Row 11
$ws1.Cells.Item(11, 1).Value2 = "Rudy Gobert"
$ws1.Cells.Item(11, 2).Value2 = 611
$ws1.Cells.Item(11, 3).Value2 = "Bryon Russell"
$ws1.Cells.Item(11, 4).Value2 = 16443
$ws1.Cells.Item(11, 5).Value2 = "Rudy Gobert"
$ws1.Cells.Item(11, 6).Value2 = 2796
$ws1.Cells.Item(11, 7).Value2 = "Deron Williams"
$ws1.Cells.Item(11, 8).Value2 = 511
$ws1.Cells.Item(11, 9).Value2 = "Deron Williams"
$ws1.Cells.Item(11, 10).Value2 = 1615
$ws1.Cells.Item(11, 11).Value2 = "Paul Millsap"
$ws1.Cells.Item(11, 12).Value2 = 3792
$ws1.Cells.Item(11, 13).Value2 = "Adrian Dantley*"
$ws1.Cells.Item(11, 14).Value2 = 1702
$ws1.Cells.Item(11, 15).Value2 = "Gordon Hayward"
$ws1.Cells.Item(11, 16).Value2 = 527
$ws1.Cells.Item(11, 17).Value2 = "Rich Kelley"
$ws1.Cells.Item(11, 18).Value2 = 498

# Column widths (character units) to match the rebuilt table layout.
# Columns 6, 7, 12, 16, 22 already have the right width from the original
# mini-tables and are intentionally left untouched.
$ws1.Columns.Item(1).ColumnWidth = 18.0
$ws1.Columns.Item(2).ColumnWidth = 12.714285714285714
$ws1.Columns.Item(3).ColumnWidth = 14.571428571428571
$ws1.Columns.Item(4).ColumnWidth = 13.428571428571429
$ws1.Columns.Item(5).ColumnWidth = 14.0
$ws1.Columns.Item(8).ColumnWidth = 12.285714285714286
$ws1.Columns.Item(9).ColumnWidth = 15.0
$ws1.Columns.Item(10).ColumnWidth = 14.571428571428571
$ws1.Columns.Item(11).ColumnWidth = 13.714285714285714
$ws1.Columns.Item(13).ColumnWidth = 11.571428571428571
$ws1.Columns.Item(15).ColumnWidth = 13.0
$ws1.Columns.Item(17).ColumnWidth = 12.0
$ws1.Columns.Item(21).ColumnWidth = 13.0

$ws1.Range("A13").Select()

# --- Sheet 2: "Portland Trail Blazers" -------------------------------------
# Only the active selection moved on this sheet; no data changed.
$ws2 = $wb.Worksheets.Item("Portland Trail Blazers")
$ws2.Range("G1:H1").Select()

$ws1.Activate()